# Update the "Förändrad" (Changed) date column (C) for every data row
# (rows 2-443) on the active worksheet from 2023-10-06 (45205) to
# 2023-10-07 (45206).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C443").Value = 45206
